$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 400
$ws.Range("I12").Value = 272.72726
$ws.Range("J12").Value = 866.6667
$ws.Range("K12").Value = 272.72726
$ws.Range("L12").Value = 866.6667
$ws.Range("M12").Value = -102.72726
$ws.Range("N12").Value = -1206.6667
$ws.Range("H15").Value = 2319.3794
$ws.Range("I15").Value = 2319.3794
$ws.Range("K15").Value = 6958.138199999999
$ws.Range("M15").Value = -6789.138199999999
$ws.Range("H32").Value = 5000
$ws.Range("J32").Value = 5000
$ws.Range("L32").Value = 5000
$ws.Range("N32").Value = -5652
$ws.Range("H38").Value = 1448.0714
$ws.Range("I38").Value = 106.166664
$ws.Range("K38").Value = 318.499992
$ws.Range("M38").Value = 53.50000799999998
$ws.Range("H40").Value = 71430240
$ws.Range("I40").Value = 1540
$ws.Range("K40").Value = 1540
$ws.Range("M40").Value = -1365
$ws.Range("H64").Value = 6399.8
$ws.Range("I64").Value = 4666.6665
$ws.Range("K64").Value = 4666.6665
$ws.Range("M64").Value = -4418.6665
$ws.Range("H67").Value = 6399.8
$ws.Range("I67").Value = 4666.6665
$ws.Range("K67").Value = 4666.6665
$ws.Range("M67").Value = -3808.6665
$ws.Range("H74").Value = 128420.22
$ws.Range("I74").Value = 143722.75
$ws.Range("K74").Value = 143722.75
$ws.Range("M74").Value = -142786.75
$ws.Range("H77").Value = 128420.22
$ws.Range("I77").Value = 143722.75
$ws.Range("K77").Value = 718613.75
$ws.Range("M77").Value = -713933.75
$ws.Range("H80").Value = 6113063.5
$ws.Range("I80").Value = 4902688.5
$ws.Range("J80").Value = 6804706.5
$ws.Range("K80").Value = 14708065.5
$ws.Range("L80").Value = 20414119.5
$ws.Range("M80").Value = -14707067.5
$ws.Range("N80").Value = -20416115.5
$ws.Range("H83").Value = 6113063.5
$ws.Range("I83").Value = 4902688.5
$ws.Range("J83").Value = 6804706.5
$ws.Range("K83").Value = 44124196.5
$ws.Range("L83").Value = 61242358.5
$ws.Range("M83").Value = -44119204.5
$ws.Range("N83").Value = -61252342.5
$ws.Range("H124").Value = 62499.5
$ws.Range("I124").Value = 25000
$ws.Range("K124").Value = 25000
$ws.Range("M124").Value = -20090
$ws.Range("H127").Value = 2173.8572
$ws.Range("I127").Value = 1702.8334
$ws.Range("K127").Value = 5108.5002
$ws.Range("M127").Value = -148.5002000000004
$ws.Range("H131").Value = 5706930.5
$ws.Range("I131").Value = 25225
$ws.Range("K131").Value = 75675
$ws.Range("M131").Value = -70635

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H102").Value = 2815.158
$ws.Range("I102").Value = 2281.8125
$ws.Range("J102").Value = 5659.6665
$ws.Range("K102").Value = 2281.8125
$ws.Range("L102").Value = 5659.6665
$ws.Range("M102").Value = -659.8125
$ws.Range("N102").Value = -8903.666499999999
$ws.Range("H122").Value = 5099.8335
$ws.Range("I122").Value = 5119.8
$ws.Range("K122").Value = 15359.4
$ws.Range("M122").Value = -12909.4
$ws.Range("H132").Value = 2946227
$ws.Range("I132").Value = 5312.5
$ws.Range("K132").Value = 15937.5
$ws.Range("M132").Value = -13407.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 8710.666999999999
$ws.Range("I20").Value = 14776.167
$ws.Range("J20").Value = 2645.1667
$ws.Range("K20").Value = 14776.167
$ws.Range("L20").Value = 2645.1667
$ws.Range("M20").Value = -14529.167
$ws.Range("N20").Value = -3139.1667
$ws.Range("H22").Value = 2323.6667
$ws.Range("I22").Value = 1813
$ws.Range("K22").Value = 1813
$ws.Range("M22").Value = -1640
$ws.Range("H105").Value = 917935.5
$ws.Range("I105").Value = 1145409.5
$ws.Range("J105").Value = 8039.8
$ws.Range("K105").Value = 1145409.5
$ws.Range("L105").Value = 8039.8
$ws.Range("M105").Value = -1143662.5
$ws.Range("N105").Value = -11533.8
$ws.Range("H107").Value = 3564.0833
$ws.Range("I107").Value = 3247.45
$ws.Range("J107").Value = 5147.25
$ws.Range("K107").Value = 3247.45
$ws.Range("L107").Value = 5147.25
$ws.Range("M107").Value = -1327.45
$ws.Range("N107").Value = -8987.25
$ws.Range("H130").Value = 29949
$ws.Range("J130").Value = 29949
$ws.Range("L130").Value = 29949
$ws.Range("N130").Value = -39989

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 33677050
$ws.Range("I31").Value = 40003976
$ws.Range("K31").Value = 40003976
$ws.Range("M31").Value = -40003681
$ws.Range("H34").Value = 33677050
$ws.Range("I34").Value = 40003976
$ws.Range("K34").Value = 40003976
$ws.Range("M34").Value = -40003774
$ws.Range("H58").Value = 2296.7058
$ws.Range("I58").Value = 2242
$ws.Range("J58").Value = 2707
$ws.Range("K58").Value = 2242
$ws.Range("L58").Value = 2707
$ws.Range("M58").Value = -2039
$ws.Range("N58").Value = -3113
$ws.Range("H94").Value = 1495.8462
$ws.Range("J94").Value = 1607.4286
$ws.Range("L94").Value = 1607.4286
$ws.Range("N94").Value = -2509.4286
$ws.Range("H105").Value = 2192.3
$ws.Range("I105").Value = 989.4286
$ws.Range("K105").Value = 989.4286
$ws.Range("M105").Value = 757.5714
$ws.Range("H136").Value = 2296.7058
$ws.Range("I136").Value = 2242
$ws.Range("J136").Value = 2707
$ws.Range("K136").Value = 6726
$ws.Range("L136").Value = 8121
$ws.Range("M136").Value = -4176
$ws.Range("N136").Value = -13221

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H26").Value = 7156.4
$ws.Range("J26").Value = 16916.5
$ws.Range("L26").Value = 50749.5
$ws.Range("N26").Value = -51325.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 927667.1
$ws.Range("I113").Value = 1560.6428
$ws.Range("J113").Value = 3088582.2
$ws.Range("K113").Value = 1560.6428
$ws.Range("L113").Value = 3088582.2
$ws.Range("M113").Value = 609.3571999999999
$ws.Range("N113").Value = -3092922.2
$ws.Range("H122").Value = 4052.5386
$ws.Range("I122").Value = 4422
$ws.Range("J122").Value = 3621.5
$ws.Range("K122").Value = 13266
$ws.Range("L122").Value = 10864.5
$ws.Range("M122").Value = -10816
$ws.Range("N122").Value = -15764.5
$ws.Range("H132").Value = 8303810
$ws.Range("I132").Value = 3588.1765
$ws.Range("K132").Value = 10764.5295
$ws.Range("M132").Value = -8234.529500000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 8203.625
$ws.Range("I61").Value = 2207
$ws.Range("K61").Value = 2207
$ws.Range("M61").Value = -2005
$ws.Range("H113").Value = 8203.625
$ws.Range("I113").Value = 2207
$ws.Range("K113").Value = 2207
$ws.Range("M113").Value = -37
$ws.Range("H122").Value = 3562.628
$ws.Range("I122").Value = 3409.6099
$ws.Range("K122").Value = 10228.8297
$ws.Range("M122").Value = -7778.8297

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H103").Value = 30166.5
$ws.Range("J103").Value = 30166.5
$ws.Range("L103").Value = 30166.5
$ws.Range("N103").Value = -32510.5
$ws.Range("H105").Value = 83472.664
$ws.Range("J105").Value = 83472.664
$ws.Range("L105").Value = 83472.664
$ws.Range("N105").Value = -90460.664
$ws.Range("H113").Value = 1014.0769
$ws.Range("I113").Value = 889.7143
$ws.Range("J113").Value = 1159.1666
$ws.Range("K113").Value = 2669.1429
$ws.Range("L113").Value = 3477.4998
$ws.Range("M113").Value = -499.1428999999998
$ws.Range("N113").Value = -7817.4998
$ws.Range("H122").Value = 4062
$ws.Range("I122").Value = 3254.3333
$ws.Range("J122").Value = 4546.6
$ws.Range("K122").Value = 9762.999899999999
$ws.Range("L122").Value = 13639.8
$ws.Range("M122").Value = -7312.999899999999
$ws.Range("N122").Value = -18539.8
$ws.Range("H132").Value = 371655.16
$ws.Range("I132").Value = 961.4761999999999
$ws.Range("K132").Value = 2884.4286
$ws.Range("M132").Value = -354.4285999999997
